$d = $word.ActiveDocument

# Locate the paragraph that ends the "Missing Persons Outlier Detection" project
# (its last bullet: "... Live Streamlit dashboard ...") so the new
# "Opportunity Intelligence Assistant" project block can be inserted immediately
# after it, and before the "AI Homelab & Active Memory Network" heading.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Live Streamlit dashboard") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph 'Live Streamlit dashboard ...'"
}

# Create a fresh empty paragraph right after the anchor, then replace that
# paragraph's content with the full new-project OOXML fragment (keeps neighboring
# paragraphs untouched, unlike InsertXML on a zero-length range sitting exactly on
# an existing paragraph boundary, which instead clobbers the adjacent paragraph).
$rng = $anchor.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$newP = $anchor.Next()
$newRng = $newP.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Opportunity Intelligence Assistant</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Agentic AI Market Analysis | Senior Living | 14 Statistical Analyses</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">• 3-call LLM pipeline directing 22 statistical methods across Census Bureau and CMS public data</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="100"/></w:pPr><w:r><w:t xml:space="preserve">• Scored Des Moines market 48.3/100; 622-word executive briefing with citation tags on every claim</w:t></w:r></w:p>'

$null = $newRng.InsertXML($xml)

Write-Output "Inserted 'Opportunity Intelligence Assistant' project block"
